# Apply "sixteenth commit with add customer fixed" edit:
# Adds two new test cases (TC_CUST_03 "Add New Customer" and TC_CUST_04 "Verify Entry")
# to the Customer_Tests worksheet, growing the used range from A1:C4 to A1:C12,
# and adjusts column widths for columns B and C on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer_Tests")
$ws.Activate()

# New rows of data (rows 5 through 12).
# NOTE: values are assigned in the same order the original author typed them
# (matching the shared-string insertion order recorded in the target workbook),
# which is not strictly top-to-bottom/left-to-right.
$ws.Range("A5").Value = "TC_CUST_03"
$ws.Range("B5").Value = "Add New Customer"

$ws.Range("C6").Value = "2.Click on ""Add New Customer"" at ""//a[contains(@class, 'btn-add')]"""

$ws.Range("C7").Value = "3.Type ""Automation User"" at ""//input[@name='name']"""

$ws.Range("C8").Value = "4.Type ""auto@nesto.com"" at ""//input[@name='email']"""

$ws.Range("C9").Value = "5.Type ""9876543210"" at ""//input[@name='mobile']"""

$ws.Range("C10").Value = "6.Click on ""//button[@type='submit']"""

$ws.Range("A11").Value = "TC_CUST_04"
$ws.Range("B11").Value = "Verify Entry"

$ws.Range("C5").Value = '1.dbexecute "{DB_QUERY}DELETE FROM customers WHERE email=''auto@nesto.com''" at ""'

$ws.Range("C11").Value = '1.Open URL "http://localhost:8080/customers"'

$ws.Range("C12").Value = "2.Verify text ""{DB_QUERY}SELECT name FROM customers WHERE email='auto@nesto.com'"" at ""//tr[td[contains(text(),'auto@nesto.com')]]/td[2]"""

# Adjust column widths on the Customer_Tests sheet (B: ~25.33, C: 128 characters)
$ws.Columns.Item(2).ColumnWidth = 24.5
$ws.Columns.Item(3).ColumnWidth = 127.16666666666667

# Update selection to the last edited cell
$ws.Range("C12").Select()
